$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bookmarks")

# New bookmark rows to append (user_id, recipe_id, created_at)
$newRows = @(
    @{ UserId = 3; RecipeId = "347"; CreatedAt = 45998.79088569444 },
    @{ UserId = 3; RecipeId = "156"; CreatedAt = 45998.790963645835 },
    @{ UserId = 3; RecipeId = "338"; CreatedAt = 45998.791139733796 },
    @{ UserId = 3; RecipeId = "204"; CreatedAt = 45998.791305023144 }
)

$startRow = 7
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $data.UserId

    # Force recipe_id to be stored as text (matches existing shared-string column)
    $ws.Range("B2").Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
    $ws.Cells.Item($row, 2).Value = $data.RecipeId

    # Reuse the existing date style from column C instead of creating a new one
    $ws.Range("C2").Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4122)
    $ws.Cells.Item($row, 3).Value = $data.CreatedAt
}

$excel.CutCopyMode = 0
